# Insert a new data row at row 512 (shifts existing rows 512:598 down to 513:599)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("512:512").Insert()

# Populate the newly inserted row 512 with the new record
$ws.Range("A512").Value = 4
$ws.Range("B512").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C512").Value = "Los Lagos"
$ws.Range("D512").Value = 45180
$ws.Range("E512").Value = 10
$ws.Range("F512").Value = 100112008
$ws.Range("G512").Value = "Coliflor"
$ws.Range("H512").Value = "Sin especificar"
$ws.Range("I512").Value = "Primera"
$ws.Range("J512").Value = 500
$ws.Range("K512").Value = 1500
$ws.Range("L512").Value = 1500
$ws.Range("M512").Value = 1500
$ws.Range("N512").Value = "`$/unidad"
$ws.Range("O512").Value = "Región Metropolitana"
$ws.Range("P512").Value = 1500
$ws.Range("Q512").Value = 1
$ws.Range("R512").Value = "Hortaliza"
